$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-140) holds the "Förändrad" (Changed) date, stored as Excel
# serial date 45205 (2023-10-06). Update it to 45206 (2023-10-07) for every
# data row, leaving all other cell contents/formatting untouched.
$lastRow = 140
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
